# Auto-applied data update mirroring the upstream scheduled-runner sync.
# Updates numeric price/profit columns (H-N) for specific rows across all 8 class sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 83337224
$ws.Range("J74").Value = 4692.3076
$ws.Range("L74").Value = 4692.3076
$ws.Range("N74").Value = -6564.3076
$ws.Range("H76").Value = 5003
$ws.Range("I76").Value = 5003
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5003
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -4688
$ws.Range("H77").Value = 83337224
$ws.Range("J77").Value = 4692.3076
$ws.Range("L77").Value = 23461.538
$ws.Range("N77").Value = -32821.538
$ws.Range("H79").Value = 5003
$ws.Range("I79").Value = 5003
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5003
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -3911
$ws.Range("H80").Value = 111112150
$ws.Range("J80").Value = 142858340
$ws.Range("L80").Value = 428575020
$ws.Range("N80").Value = -428577016
$ws.Range("H83").Value = 111112150
$ws.Range("J83").Value = 142858340
$ws.Range("L83").Value = 1285725060
$ws.Range("N83").Value = -1285735044
$ws.Range("H112").Value = 2798.5186
$ws.Range("J112").Value = 2890.82
$ws.Range("L112").Value = 8672.460000000001
$ws.Range("N112").Value = -10888.46
$ws.Range("H116").Value = 31868104
$ws.Range("I116").Value = 49247396
$ws.Range("J116").Value = 6065.6665
$ws.Range("K116").Value = 49247396
$ws.Range("L116").Value = 6065.6665
$ws.Range("M116").Value = -49243954
$ws.Range("N116").Value = -12949.6665
$ws.Range("H121").Value = 4728.2666
$ws.Range("J121").Value = 4728.2666
$ws.Range("L121").Value = 14184.7998
$ws.Range("N121").Value = -17678.7998
$ws.Range("H138").Value = 3223.37
$ws.Range("I138").Value = 1422.48
$ws.Range("J138").Value = 3823.6667
$ws.Range("K138").Value = 4267.440000000001
$ws.Range("L138").Value = 11471.0001
$ws.Range("M138").Value = 872.5599999999995
$ws.Range("N138").Value = -21751.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5254.6055
$ws.Range("I61").Value = 5102.0835
$ws.Range("K61").Value = 5102.0835
$ws.Range("M61").Value = -4890.0835
$ws.Range("H74").Value = 1094
$ws.Range("I74").Value = 825.9474
$ws.Range("J74").Value = 1659.8889
$ws.Range("K74").Value = 825.9474
$ws.Range("L74").Value = 1659.8889
$ws.Range("M74").Value = 48.05259999999998
$ws.Range("N74").Value = -3407.8889
$ws.Range("H77").Value = 1094
$ws.Range("I77").Value = 825.9474
$ws.Range("J77").Value = 1659.8889
$ws.Range("K77").Value = 4129.737
$ws.Range("L77").Value = 8299.4445
$ws.Range("M77").Value = 238.2629999999999
$ws.Range("N77").Value = -17035.4445
$ws.Range("H122").Value = 2662.0952
$ws.Range("I122").Value = 1429.3226
$ws.Range("K122").Value = 4287.9678
$ws.Range("M122").Value = -1837.9678
$ws.Range("H136").Value = 5254.6055
$ws.Range("I136").Value = 5102.0835
$ws.Range("K136").Value = 15306.2505
$ws.Range("M136").Value = -12756.2505

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1604562.4
$ws.Range("I99").Value = 2606039
$ws.Range("J99").Value = 2200
$ws.Range("K99").Value = 2606039
$ws.Range("L99").Value = 2200
$ws.Range("M99").Value = -2604541
$ws.Range("N99").Value = -5196
$ws.Range("H105").Value = 1498.0667
$ws.Range("I105").Value = 1498.0667
$ws.Range("K105").Value = 1498.0667
$ws.Range("M105").Value = 248.9332999999999
$ws.Range("H134").Value = 1306.3529
$ws.Range("I134").Value = 1270.7241
$ws.Range("J134").Value = 1513
$ws.Range("K134").Value = 3812.1723
$ws.Range("L134").Value = 4539
$ws.Range("M134").Value = -1277.1723
$ws.Range("N134").Value = -9609

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1055.25
$ws.Range("I22").Value = 481.4
$ws.Range("J22").Value = 2011.6666
$ws.Range("K22").Value = 481.4
$ws.Range("L22").Value = 2011.6666
$ws.Range("M22").Value = -131.4
$ws.Range("N22").Value = -2711.6666
$ws.Range("H31").Value = 3838.9673
$ws.Range("I31").Value = 1243.5834
$ws.Range("K31").Value = 1243.5834
$ws.Range("M31").Value = -948.5834
$ws.Range("H34").Value = 3838.9673
$ws.Range("I34").Value = 1243.5834
$ws.Range("K34").Value = 1243.5834
$ws.Range("M34").Value = -1041.5834
$ws.Range("H132").Value = 6811872
$ws.Range("I132").Value = 7415002.5
$ws.Range("K132").Value = 22245007.5
$ws.Range("M132").Value = -22242477.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1213
$ws.Range("I5").Value = 1213
$ws.Range("K5").Value = 3639
$ws.Range("M5").Value = -3527
$ws.Range("H118").Value = 2680.5557
$ws.Range("I118").Value = 1870.8334
$ws.Range("J118").Value = 4300
$ws.Range("K118").Value = 5612.5002
$ws.Range("L118").Value = 12900
$ws.Range("M118").Value = -4369.5002
$ws.Range("N118").Value = -15386
$ws.Range("H132").Value = 4392.353
$ws.Range("I132").Value = 893.1667
$ws.Range("J132").Value = 6301
$ws.Range("K132").Value = 8038.5003
$ws.Range("L132").Value = 56709
$ws.Range("M132").Value = -5508.5003
$ws.Range("N132").Value = -61769
$ws.Range("H135").Value = 1213
$ws.Range("I135").Value = 1213
$ws.Range("K135").Value = 10917
$ws.Range("M135").Value = -8382

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("H70").Value = 7581049
$ws.Range("I70").Value = 22729148
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 22729148
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -22728878
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 7581049
$ws.Range("I73").Value = 22729148
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 22729148
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -22728212
$ws.Range("N73").Value = -8872
$ws.Range("H113").Value = 4374.75
$ws.Range("I113").Value = 3999.5
$ws.Range("J113").Value = 4750
$ws.Range("K113").Value = 3999.5
$ws.Range("L113").Value = 4750
$ws.Range("M113").Value = -1829.5
$ws.Range("N113").Value = -9090

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 7200
$ws.Range("J14").Value = 7200
$ws.Range("L14").Value = 7200
$ws.Range("N14").Value = -7544
$ws.Range("H22").Value = 1046.9375
$ws.Range("I22").Value = 925.125
$ws.Range("J22").Value = 1168.75
$ws.Range("K22").Value = 925.125
$ws.Range("L22").Value = 1168.75
$ws.Range("M22").Value = -630.125
$ws.Range("N22").Value = -1758.75
$ws.Range("H27").Value = 1046.9375
$ws.Range("I27").Value = 925.125
$ws.Range("J27").Value = 1168.75
$ws.Range("K27").Value = 925.125
$ws.Range("L27").Value = 1168.75
$ws.Range("M27").Value = -818.125
$ws.Range("N27").Value = -1382.75
$ws.Range("H132").Value = 3712.2144
$ws.Range("I132").Value = 3576.457
$ws.Range("K132").Value = 10729.371
$ws.Range("M132").Value = -8199.370999999999
$ws.Range("H136").Value = 2596.3235
$ws.Range("I136").Value = 1828.7407
$ws.Range("K136").Value = 5486.2221
$ws.Range("M136").Value = -2936.2221

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 10999.75
$ws.Range("J11").Value = 10999.75
$ws.Range("L11").Value = 10999.75
$ws.Range("N11").Value = -11283.75
$ws.Range("H13").Value = 2899.3333
$ws.Range("J13").Value = 2899.3333
$ws.Range("L13").Value = 2899.3333
$ws.Range("N13").Value = -3179.3333
$ws.Range("H17").Value = 900
$ws.Range("I17").Value = 900
$ws.Range("K17").Value = 900
$ws.Range("M17").Value = -728
$ws.Range("H49").Value = 53999
$ws.Range("J49").Value = 53999
$ws.Range("L49").Value = 53999
$ws.Range("N49").Value = -54459
$ws.Range("H107").Value = 14519
$ws.Range("I107").Value = 15149
$ws.Range("J107").Value = 11999
$ws.Range("K107").Value = 45447
$ws.Range("L107").Value = 35997
$ws.Range("M107").Value = -43527
$ws.Range("N107").Value = -39837
$ws.Range("H132").Value = 13890196
$ws.Range("I132").Value = 1353.8334
$ws.Range("J132").Value = 83334410
$ws.Range("K132").Value = 4061.5002
$ws.Range("L132").Value = 250003230
$ws.Range("M132").Value = -1531.5002
$ws.Range("N132").Value = -250008290
